$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Reverted blocks_final file to default": restore the editable trial-count
# columns (C:F) from the hidden default columns (M:P) for each data row.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 13).Value2  # C = M
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 14).Value2  # D = N
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 15).Value2  # E = O
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 16).Value2  # F = P
}

# Update the selected cell to match the saved view state.
$ws.Range("F14").Select()
